# Contacts list cleanup: remove stale contact rows and refresh remaining
# entries with the current Phone / Name / Group values (used by the new
# global "search for all modules" feature).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had 4 contact rows; only the first two are still valid.
# Drop rows 3 and 4 entirely.
$ws.Rows("3:4").Delete()

# Refresh the surviving rows with the current contact data.
$ws.Cells.Item(1, 1).Value = "09176214704"
$ws.Cells.Item(1, 2).Value = "Jasper Barcelona"
$ws.Cells.Item(1, 3).Value = "New Group"

$ws.Cells.Item(2, 1).Value = "09772312533"
$ws.Cells.Item(2, 2).Value = "Leanza Etorma"
$ws.Cells.Item(2, 3).Value = "New Group"

# Let the Name/Phone columns re-fit their (now shorter) contents.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

# Matches the selection left behind in the saved workbook.
$null = $ws.Range("A10").Select()
